$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1780415430267062
$ws.Range("C2").Value = 0.5905044510385756
$ws.Range("J2").Value = 0.01780415430267062
$ws.Range("P2").Value = 0.142433234421365
$ws.Range("S2").Value = 0.0712166172106825
# Row 3
$ws.Range("B3").Value = 0.004975124378109453
$ws.Range("C3").Value = 0.004975124378109453
$ws.Range("J3").Value = 0.03482587064676617
$ws.Range("P3").Value = 0.8109452736318408
$ws.Range("S3").Value = 0.1442786069651741
# Row 4
$ws.Range("J4").Value = 0.07017543859649122
$ws.Range("P4").Value = 0.8070175438596491
$ws.Range("S4").Value = 0.1228070175438596
# Row 6
$ws.Range("B6").Value = 0.06010928961748634
$ws.Range("D6").Value = 0.03278688524590164
$ws.Range("E6").Value = 0.00546448087431694
$ws.Range("F6").Value = 0.04371584699453552
$ws.Range("J6").Value = 0.180327868852459
$ws.Range("O6").Value = 0.03278688524590164
$ws.Range("Q6").Value = 0.2513661202185792
$ws.Range("R6").Value = 0.06557377049180328
$ws.Range("S6").Value = 0.3278688524590164
# Row 7
$ws.Range("B7").Value = 0.109452736318408
$ws.Range("D7").Value = 0.03482587064676617
$ws.Range("E7").Value = 0.004975124378109453
$ws.Range("F7").Value = 0.06965174129353234
$ws.Range("J7").Value = 0.154228855721393
$ws.Range("O7").Value = 0.01492537313432836
$ws.Range("Q7").Value = 0.1990049751243781
$ws.Range("R7").Value = 0.0845771144278607
$ws.Range("S7").Value = 0.3283582089552239
# Row 8
$ws.Range("B8").Value = 0.1129411764705882
$ws.Range("D8").Value = 0.03058823529411765
$ws.Range("E8").Value = 0.002352941176470588
$ws.Range("F8").Value = 0.04
$ws.Range("J8").Value = 0.1105882352941177
$ws.Range("O8").Value = 0.02588235294117647
$ws.Range("Q8").Value = 0.2517647058823529
$ws.Range("R8").Value = 0.04705882352941176
$ws.Range("S8").Value = 0.3788235294117647
# Row 9
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.02651515151515152
$ws.Range("E9").Value = 0.003787878787878788
$ws.Range("F9").Value = 0.02651515151515152
$ws.Range("J9").Value = 0.1401515151515151
$ws.Range("O9").Value = 0.02272727272727273
$ws.Range("Q9").Value = 0.2348484848484849
$ws.Range("R9").Value = 0.04924242424242424
$ws.Range("S9").Value = 0.4053030303030303
# Row 10
$ws.Range("B10").Value = 0.1391585760517799
$ws.Range("D10").Value = 0.02184466019417476
$ws.Range("E10").Value = 0.002427184466019417
$ws.Range("F10").Value = 0.07605177993527508
$ws.Range("J10").Value = 0.1124595469255663
$ws.Range("O10").Value = 0.007281553398058253
$ws.Range("Q10").Value = 0.2192556634304207
$ws.Range("R10").Value = 0.07281553398058252
$ws.Range("S10").Value = 0.348705501618123
# Row 11
$ws.Range("G11").Value = 0.1211180124223603
$ws.Range("J11").Value = 0.08695652173913043
$ws.Range("K11").Value = 0.1832298136645963
$ws.Range("L11").Value = 0.5838509316770186
$ws.Range("S11").Value = 0.02484472049689441
# Row 12
$ws.Range("G12").Value = 0.6979166666666666
$ws.Range("J12").Value = 0.2447916666666667
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.02604166666666667
$ws.Range("S12").Value = 0.02604166666666667
# Row 13
$ws.Range("G13").Value = 0.6382978723404256
$ws.Range("J13").Value = 0.3191489361702128
$ws.Range("S13").Value = 0.0425531914893617
# Row 15
$ws.Range("F15").Value = 0.005208333333333333
$ws.Range("H15").Value = 0.1145833333333333
$ws.Range("I15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.3541666666666667
$ws.Range("K15").Value = 0.109375
$ws.Range("M15").Value = 0.015625
$ws.Range("O15").Value = 0.03645833333333334
$ws.Range("S15").Value = 0.28125
# Row 16
$ws.Range("F16").Value = 0.007604562737642586
$ws.Range("H16").Value = 0.2015209125475285
$ws.Range("I16").Value = 0.09885931558935361
$ws.Range("J16").Value = 0.3307984790874525
$ws.Range("K16").Value = 0.1216730038022814
$ws.Range("M16").Value = 0.02281368821292776
$ws.Range("N16").Value = 0.007604562737642586
$ws.Range("O16").Value = 0.05703422053231939
$ws.Range("S16").Value = 0.1520912547528517
# Row 17
$ws.Range("F17").Value = 0.007692307692307693
$ws.Range("H17").Value = 0.1423076923076923
$ws.Range("I17").Value = 0.1115384615384615
$ws.Range("J17").Value = 0.425
$ws.Range("K17").Value = 0.1192307692307692
$ws.Range("M17").Value = 0.01730769230769231
$ws.Range("O17").Value = 0.04038461538461539
$ws.Range("S17").Value = 0.1365384615384615
# Row 18
$ws.Range("F18").Value = 0.01333333333333333
$ws.Range("H18").Value = 0.1866666666666667
$ws.Range("I18").Value = 0.1333333333333333
$ws.Range("J18").Value = 0.3666666666666666
$ws.Range("K18").Value = 0.08
$ws.Range("M18").Value = 0.02
$ws.Range("O18").Value = 0.07333333333333333
$ws.Range("S18").Value = 0.1266666666666667
# Row 19
$ws.Range("F19").Value = 0.009700889248181084
$ws.Range("H19").Value = 0.2029102667744543
$ws.Range("I19").Value = 0.1156022635408246
$ws.Range("J19").Value = 0.3459983831851253
$ws.Range("K19").Value = 0.1050929668552951
$ws.Range("M19").Value = 0.02425222312045271
$ws.Range("N19").Value = 0.001616814874696847
$ws.Range("O19").Value = 0.05901374292643492
$ws.Range("S19").Value = 0.1358124494745352

Write-Output "Applied 113 cell updates to transition matrix"
